$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "66.649.90"
$ws.Cells.Item(2, 5).Value = "  -0.74%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.087.56"
$ws.Cells.Item(3, 5).Value = "  -1.50%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.14%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "575.97"
$ws.Cells.Item(5, 5).Value = "  -0.93%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "172.10"
$ws.Cells.Item(6, 5).Value = "  -1.18%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 5).Value = "  +0.06%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.084.93"
$ws.Cells.Item(8, 5).Value = "  -1.35%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.512"
$ws.Cells.Item(9, 5).Value = "  -2.12%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "6.34"
$ws.Cells.Item(10, 5).Value = "  -1.60%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.151"
$ws.Cells.Item(11, 5).Value = "  -2.89%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.469"
$ws.Cells.Item(12, 5).Value = "  -2.59%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000238"
$ws.Cells.Item(13, 5).Value = "  -4.43%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "36.01"
$ws.Cells.Item(14, 5).Value = "  -4.39%  "
$ws.Cells.Item(15, 5).Value = "  -1.03%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.604.66"
$ws.Cells.Item(16, 5).Value = "  -1.16%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "66.623.78"
$ws.Cells.Item(17, 5).Value = "  -0.64%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "6.95"
$ws.Cells.Item(18, 5).Value = "  -2.82%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "16.76"
$ws.Cells.Item(19, 5).Value = "  +2.12%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "3.087.61"
$ws.Cells.Item(20, 5).Value = "  -1.43%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "486.88"
$ws.Cells.Item(21, 5).Value = "  -1.24%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "7.81"
$ws.Cells.Item(22, 5).Value = "  -1.14%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.686"
$ws.Cells.Item(23, 5).Value = "  -3.49%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "83.34"
$ws.Cells.Item(24, 5).Value = "  -1.17%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "12.68"
$ws.Cells.Item(25, 5).Value = "  -4.55%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.22"
$ws.Cells.Item(26, 5).Value = "  -3.66%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.02"
$ws.Cells.Item(27, 5).Value = "  -3.48%  "
$ws.Cells.Item(28, 5).Value = "  -0.03%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.96"
$ws.Cells.Item(29, 5).Value = "  -0.16%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.24"
$ws.Cells.Item(30, 5).Value = "  -5.29%  "
$ws.Cells.Item(31, 5).Value = "  -3.98%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "27.97"
$ws.Cells.Item(32, 5).Value = "  -2.95%  "
$ws.Cells.Item(33, 5).Value = "  -3.11%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0₃0928"
$ws.Cells.Item(34, 5).Value = "  -2.93%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.999"
$ws.Cells.Item(35, 5).Value = "  +0.08%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "48.99"
$ws.Cells.Item(36, 5).Value = "  +4.27%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.58"
$ws.Cells.Item(37, 5).Value = "  -5.73%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.938"
$ws.Cells.Item(38, 5).Value = "  -4.27%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "49.07"
$ws.Cells.Item(39, 5).Value = "  -2.09%  "
$ws.Cells.Item(40, 5).Value = "  -2.07%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.122"
$ws.Cells.Item(41, 5).Value = "  -1.44%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.96"
$ws.Cells.Item(42, 5).Value = "  -5.08%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "8.23"
$ws.Cells.Item(43, 5).Value = "  -4.06%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.58"
$ws.Cells.Item(44, 5).Value = "  -1.50%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.774.85"
$ws.Cells.Item(45, 5).Value = "  -2.14%  "
$ws.Cells.Item(46, 5).Value = "  -3.10%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "367.49"
$ws.Cells.Item(47, 5).Value = "  -5.11%  "
$ws.Cells.Item(48, 5).Value = "  -1.32%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "24.36"
$ws.Cells.Item(50, 5).Value = "  -3.09%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.17"
$ws.Cells.Item(51, 5).Value = "  -2.75%  "
